$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 31.01.2022 01:30"

# D5: was text "+0.4" -> now numeric 0.4
$ws.Range("D5").Value = 0.4

# E5: was text "2022-01-31 01:15:09" -> now numeric date serial with date/time format
$ws.Range("E5").Value = 44592.0521875
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
